# Adding the changes we made on may 9th
#
# The accelerometer log gained 11 new readings that belong at the top of
# the table (rows 2-12), pushing every existing reading down by 11 rows.
# The table also drops its previous last row, since the new bottom of the
# data (row 31) matches the pre-existing second-to-last reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing readings (rows 2:21) down by 11 rows, making room for
# the 11 freshly recorded samples at the top of the table.
$ws.Rows("2:12").Insert()

# The insert inherits the bold header formatting from row 1 - the source
# data rows carry no explicit styling, so reset the newly inserted rows
# back to the default/unstyled look.
$ws.Range("A2:C12").Style = "Normal"

# New accelerometer samples recorded on may 9th, inserted at the top.
$newSamples = @(
    @(-1.610764980316162, 1.577288150787354, 0.315173327922821),
    @(-1.949196338653564, 1.61798095703125, 0.4421060681343078),
    @(-1.840867042541504, 1.58759355545044, 0.5192338824272156),
    @(-1.788941383361816, 1.582527160644531, 0.4804926216602325),
    @(-1.827802658081055, 1.591060638427734, 0.4033206701278686),
    @(-1.872776985168457, 1.712420463562011, 0.4200127720832824),
    @(-1.756282329559326, 1.523788452148438, 0.3277221620082855),
    @(-1.891244411468506, 1.559478759765625, 0.1866782307624817),
    @(-1.622483730316162, 1.571603775024414, 0.2709611356258392),
    @(-1.867420673370361, 1.565328121185303, 0.3432579040527344),
    @(-1.858330726623535, 1.508580207824707, 0.2882210314273834)
)

$row = 2
foreach ($sample in $newSamples) {
    $ws.Cells.Item($row, 1).Value = $sample[0]
    $ws.Cells.Item($row, 2).Value = $sample[1]
    $ws.Cells.Item($row, 3).Value = $sample[2]
    $row = $row + 1
}

# The insert pushed the original last reading (previously row 21) down to
# row 32, beyond the new end of the table (row 31) - drop it.
$ws.Rows("32:32").Delete()
